$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before A, shifting the whole table one column
# to the right (old A -> B, old B -> C, ... old F -> G).
$ws.Columns("A").Insert()

# --- Header row -----------------------------------------------------------
# Old A1 was empty, so after the shift B1 is empty too. Give it the new
# "segments" header, matching the style already used by the other headers.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "segments"

# --- New index column (A) --------------------------------------------------
# Copy the header-like style (bold, bordered, centered) that used to live on
# the old column A (now on column B) onto the new column A.
$ws.Range("B2").Copy()
$ws.Range("A2:A20").PasteSpecial(-4122)

# 19 category rows (rows 2-20) get a 0-based numeric index in column A.
for ($i = 0; $i -le 18; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $i
}

# --- Strip the old header-like style off column B (now plain text) --------
$ws.Range("C2").Copy()
$ws.Range("B2:B20").PasteSpecial(-4122)
